# Update lipid/recruitment data prep loop results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New meanLipids (column B) values for rows 2-12 (years 2011-2021)
$meanLipids = @(
    17.041290322580643,
    0,
    5.6309756097560975,
    12.558092105263157,
    11.038570620938627,
    8.9742483660130716,
    2.986470588235294,
    7.3007031914624276,
    17.456099290780141,
    30.863691189496603,
    17.204920088467741
)

# New recruitment (column C) values for rows 2-10 (years 2011-2019); rows 11-12 now have no value
$recruitment = @(
    0.83840000000000003,
    0.87959999999999994,
    0.74650000000000005,
    0.4955,
    0.41660000000000003,
    0.30133333333333329,
    0.70124999999999993,
    0.25166666666666665,
    0.69133333333333324
)

for ($i = 0; $i -lt $meanLipids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $meanLipids[$i]
}

for ($i = 0; $i -lt $recruitment.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $recruitment[$i]
}

# Remove the now-empty recruitment cells for 2020 (row 11) and 2021 (row 12)
$ws.Cells.Item(11, 3).ClearContents()
$ws.Cells.Item(12, 3).ClearContents()

# Resize column B to reflect the new (bestFit) width captured in the diff.
# (Target stored width is ~11.73; this is the closest value this runtime's
# column-width quantization can reach.)
$ws.Columns.Item(2).ColumnWidth = 10.83

# Update the active selection to match the authored state
$ws.Range("C10").Select() | Out-Null

$wb.Save()
